$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 287, pushing the existing rows (old 287-313)
# down to 289-315.
$ws.Rows("287:288").Insert()

# New row 287 data
$ws.Cells.Item(287, 1).Value = 5
$ws.Cells.Item(287, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(287, 3).Value = "Maule"
$ws.Cells.Item(287, 4).Value = 44783
$ws.Cells.Item(287, 5).Value = 7
$ws.Cells.Item(287, 6).Value = "Fruta"
$ws.Cells.Item(287, 7).Value = 100101
$ws.Cells.Item(287, 8).Value = "Berries"
$ws.Cells.Item(287, 9).Value = 100101007
$ws.Cells.Item(287, 10).Value = "Kiwi"
$ws.Cells.Item(287, 11).Value = "Hayward"
$ws.Cells.Item(287, 12).Value = "Primera"
$ws.Cells.Item(287, 13).Value = 260
$ws.Cells.Item(287, 14).Value = 7000
$ws.Cells.Item(287, 15).Value = 7000
$ws.Cells.Item(287, 16).Value = 7000
$ws.Cells.Item(287, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(287, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(287, 19).Value = 389
$ws.Cells.Item(287, 20).Value = 18

# New row 288 data
$ws.Cells.Item(288, 1).Value = 5
$ws.Cells.Item(288, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(288, 3).Value = "Maule"
$ws.Cells.Item(288, 4).Value = 44783
$ws.Cells.Item(288, 5).Value = 7
$ws.Cells.Item(288, 6).Value = "Fruta"
$ws.Cells.Item(288, 7).Value = 100101
$ws.Cells.Item(288, 8).Value = "Berries"
$ws.Cells.Item(288, 9).Value = 100101007
$ws.Cells.Item(288, 10).Value = "Kiwi"
$ws.Cells.Item(288, 11).Value = "Hayward"
$ws.Cells.Item(288, 12).Value = "Segunda"
$ws.Cells.Item(288, 13).Value = 200
$ws.Cells.Item(288, 14).Value = 5000
$ws.Cells.Item(288, 15).Value = 5000
$ws.Cells.Item(288, 16).Value = 5000
$ws.Cells.Item(288, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(288, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(288, 19).Value = 278
$ws.Cells.Item(288, 20).Value = 18
